# QC sprint 1 - first stage
# Teacher.xlsx template: remove the Datebirth / Gender / Nationaid columns
# (these were columns C, D and E) so the sheet goes from
# Firstname, Lastname, Datebirth, Gender, Nationaid, Email, Password, Phone, Username
# down to
# Firstname, Lastname, Email, Password, Phone, Username

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the three unwanted columns (C:E). Excel shifts everything
# after them to the left, which also drops the now-unused shared
# strings ("Datebirth", "Gender", "Nationaid") and shrinks the used
# dimension from A1:I11 down to A1:F11.
$ws.Columns("C:E").Delete()

# Tidy up the remaining column widths (the surviving columns keep
# roughly their old widths, but nudge them to the refreshed layout).
$ws.Columns.Item(1).ColumnWidth = 31.833333333333332
$ws.Columns.Item(2).ColumnWidth = 28.666666666666668
$ws.Columns.Item(3).ColumnWidth = 33.666666666666664
$ws.Columns.Item(4).ColumnWidth = 22.166666666666668
$ws.Columns.Item(5).ColumnWidth = 21.666666666666668
$ws.Columns.Item(6).ColumnWidth = 22.666666666666668

# Normalize the (thick-border) row heights for the header + the 10
# data rows.
for ($r = 1; $r -le 11; $r++) {
    $ws.Rows.Item($r).RowHeight = 15
}

# Leave the cursor where the author left it before saving.
$null = $ws.Range("E6").Select()

$wb.Save()
